# Refresh the cryptos list with the latest coinranking.com snapshot:
# updated prices / 1h volume % for (almost) every row, plus two pairs
# of rows that swapped rank order (WrappedEther/TRON, Quant/TheSandbox).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text such as "30.963.11" or "0.4881".
# A leading apostrophe forces Excel's COM layer to keep the literal text
# instead of re-parsing it as a number (which would e.g. turn "1.000"
# into 1, or render "0.000007689" in scientific notation).

$ws.Range('D2').Value = '''31.033.63'
$ws.Range('E2').Value = '  +1.16%  '

$ws.Range('D3').Value = '''1.954.28'
$ws.Range('E3').Value = '  -0.40%  '

$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '''245.36'
$ws.Range('E5').Value = '  -1.39%  '

$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').Value = '''0.4883'
$ws.Range('E7').Value = '  +0.93%  '

$ws.Range('D8').Value = '''0.2949'
$ws.Range('E8').Value = '  -0.20%  '

$ws.Range('D9').Value = '''0.06829'
$ws.Range('E9').Value = '  +0.57%  '

$ws.Range('D10').Value = '''19.20'
$ws.Range('E10').Value = '  -0.76%  '

$ws.Range('D11').Value = '''107.15'
$ws.Range('E11').Value = '  -3.42%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.07799'
$ws.Range('E12').Value = '  +0.73%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '''1.940.44'
$ws.Range('E13').Value = '  -1.23%  '

$ws.Range('D14').Value = '''5.461'
$ws.Range('E14').Value = '  -0.59%  '

$ws.Range('D15').Value = '''0.7019'
$ws.Range('E15').Value = '  +1.52%  '

$ws.Range('D16').Value = '''283.07'
$ws.Range('E16').Value = '  -3.73%  '

$ws.Range('D17').Value = '''31.064.44'
$ws.Range('E17').Value = '  +1.21%  '

$ws.Range('D18').Value = '''13.20'
$ws.Range('E18').Value = '  -0.85%  '

$ws.Range('D19').Value = '''0.000007689'
$ws.Range('E19').Value = '  -0.04%  '

$ws.Range('D20').Value = '''2.206.47'
$ws.Range('E20').Value = '  -0.38%  '

$ws.Range('E21').Value = '  -0.01%  '

$ws.Range('D22').Value = '''5.479'
$ws.Range('E22').Value = '  -3.43%  '

$ws.Range('D23').Value = '''1.000'
$ws.Range('E23').Value = '  +0.06%  '

$ws.Range('D24').Value = '''6.488'
$ws.Range('E24').Value = '  -1.87%  '

$ws.Range('D25').Value = '''9.829'
$ws.Range('E25').Value = '  -1.07%  '

$ws.Range('D26').Value = '''169.53'
$ws.Range('E26').Value = '  -0.74%  '

$ws.Range('D27').Value = '''20.02'
$ws.Range('E27').Value = '  -0.71%  '

$ws.Range('D28').Value = '''2.203'
$ws.Range('E28').Value = '  +0.13%  '

$ws.Range('D29').Value = '''0.1057'
$ws.Range('E29').Value = '  -1.51%  '

$ws.Range('D30').Value = '''1.417'
$ws.Range('E30').Value = '  -1.78%  '

$ws.Range('D31').Value = '''1.579'
$ws.Range('E31').Value = '  -1.60%  '

$ws.Range('D32').Value = '''4.595'
$ws.Range('E32').Value = '  -2.14%  '

$ws.Range('D33').Value = '''4.450'
$ws.Range('E33').Value = '  -0.48%  '

$ws.Range('D34').Value = '''0.04954'
$ws.Range('E34').Value = '  -3.53%  '

$ws.Range('D35').Value = '''0.7671'
$ws.Range('E35').Value = '  -1.64%  '

$ws.Range('D36').Value = '''1.172'
$ws.Range('E36').Value = '  -0.57%  '

$ws.Range('D37').Value = '''2.729'
$ws.Range('E37').Value = '  -0.25%  '

$ws.Range('D38').Value = '''0.02008'
$ws.Range('E38').Value = '  -2.61%  '

$ws.Range('D39').Value = '''2.704'
$ws.Range('E39').Value = '  -0.30%  '

$ws.Range('D40').Value = '''6.503'
$ws.Range('E40').Value = '  +6.00%  '

$ws.Range('D41').Value = '''2.127'
$ws.Range('E41').Value = '  +2.61%  '

$ws.Range('D42').Value = '''74.30'
$ws.Range('E42').Value = '  +5.65%  '

$ws.Range('D43').Value = '''0.8855'
$ws.Range('E43').Value = '  +1.18%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.4461'
$ws.Range('E44').Value = '  -0.15%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '''109.33'
$ws.Range('E45').Value = '  -1.73%  '

$ws.Range('D46').Value = '''8.152'
$ws.Range('E46').Value = '  +10.25%  '

$ws.Range('E47').Value = '  -0.06%  '

$ws.Range('D48').Value = '''992.53'
$ws.Range('E48').Value = '  +9.68%  '

$ws.Range('D49').Value = '''0.1262'
$ws.Range('E49').Value = '  -1.45%  '

$ws.Range('D50').Value = '''9.335'
$ws.Range('E50').Value = '  -0.88%  '

$ws.Range('D51').Value = '''0.2597'
$ws.Range('E51').Value = '  +3.17%  '
